$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title meta tag: `1 What is Blender ` -> ` The Shrink and Fatten Tool `
#    (content ends up as: content=" The Shrink and Fatten Tool ">)
# ---------------------------------------------------------------------------
$rTitle = $d.Content
$null = $rTitle.Find.Execute("1 What is Blender")
$rTitle.Text = " The Shrink and Fatten Tool"

# ---------------------------------------------------------------------------
# 2) Keywords meta tag: insert "The Shrink and Fatten Tool, " before the
#    "Blender, 3D Modeling, Animation, Graphic Art" list.
#    (This is the FIRST occurrence of that phrase in the document.)
# ---------------------------------------------------------------------------
$rKeywords = $d.Content
$null = $rKeywords.Find.Execute("Blender, 3D Modeling, Animation, Graphic Art")
$insKeywords = $d.Range($rKeywords.Start, $rKeywords.Start)
$insKeywords.InsertBefore("The Shrink and Fatten Tool, ")

# ---------------------------------------------------------------------------
# 3) Description meta tag: replace
#       what the 3D modeling program "Blender " is all about./
#    with
#       about The Shrink and Fatten Tool, in Edit mode, while working in the Blender Application/
# ---------------------------------------------------------------------------
$rDescStart = $d.Content
$null = $rDescStart.Find.Execute("what the 3D modeling program")
$rDescEnd = $d.Content
$null = $rDescEnd.Find.Execute("about./>")
$rDesc = $d.Range($rDescStart.Start, $rDescEnd.End)
$rDesc.Text = "about The Shrink and Fatten Tool, in Edit mode, while working in the Blender Application/>"

# ---------------------------------------------------------------------------
# 4) Category meta tag: insert "The Shrink and Fatten Tool, " before the
#    "Blender, 3D Modeling, Animation, Graphic Art" list.
#    Scope the search to start after the "category" label so that the
#    keywords occurrence earlier in the document (which, after step 2,
#    also now contains this phrase) is not matched again.
# ---------------------------------------------------------------------------
$rCategoryLabel = $d.Content
$null = $rCategoryLabel.Find.Execute("category")
$rCategory = $d.Range($rCategoryLabel.End, $d.Content.End)
$null = $rCategory.Find.Execute("Blender, 3D Modeling, Animation, Graphic Art")
$insCategory = $d.Range($rCategory.Start, $rCategory.Start)
$insCategory.InsertBefore("The Shrink and Fatten Tool, ")

# ---------------------------------------------------------------------------
# 5) Revised date
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("Wednesday, December 11, 2024", $true, $false, $false, $false, $false, $true, 1, $false, "Saturday, January 25, 2025", 2)

# ---------------------------------------------------------------------------
# 6) URL path
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("Enlightenment/Articles/2024/8-Blender-2024/1-What-Is-Blender/1-What-Is-Blender.html", $true, $false, $false, $false, $false, $true, 1, $false, "Enlightenment/Articles/2025/1-Blender-Continued/2-Edit-Mode/1-The-Menus/1-The-Tools-Menu/11-The-Shrink-And-Fatten-Tool/The-Shrink-And-Fatten-Tool.html", 2)
